$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
Write-Host "Name: $($sh.Name)"
$tf = $sh.TextFrame
$tr = $tf.TextRange
Write-Host "Text: [$($tr.Text)]"
$tr.Text = "23-mar-18"
